$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, shifting the existing rows
# 116-154 down to 117-155 (dimension grows from R154 to R155).
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new weekly price record.
$ws.Cells.Item(116, 1).Value = 7
$ws.Cells.Item(116, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(116, 3).Value = "Ñuble"
$ws.Cells.Item(116, 4).Value = 44524
$ws.Cells.Item(116, 5).Value = 16
$ws.Cells.Item(116, 6).Value = 100112006
$ws.Cells.Item(116, 7).Value = "Repollo"
$ws.Cells.Item(116, 8).Value = "Crespo record"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 240
$ws.Cells.Item(116, 11).Value = 600
$ws.Cells.Item(116, 12).Value = 700
$ws.Cells.Item(116, 13).Value = 650
$ws.Cells.Item(116, 14).Value = "`$/unidad"
$ws.Cells.Item(116, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(116, 16).Value = 650
$ws.Cells.Item(116, 17).Value = 1
$ws.Cells.Item(116, 18).Value = "Hortaliza"
